$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '69.051.71'
$ws.Range("E2").Value = '  -2.16%  '
$ws.Range("D3").Value = '2.509.97'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue "D5" '571.51'
$ws.Range("E5").Value = '  -0.80%  '
Set-TextValue "D6" '166.28'
$ws.Range("E6").Value = '  -1.92%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").Value = '2.510.87'
$ws.Range("E9").Value = '  -0.75%  '
Set-TextValue "D10" '0.160'
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("E11").Value = '  -0.52%  '
Set-TextValue "D12" '0.352'
$ws.Range("E12").Value = '  +2.75%  '
Set-TextValue "D13" '4.91'
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("D14").Value = '2.974.78'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = '69.119.23'
$ws.Range("E15").Value = '  -1.85%  '
Set-TextValue "D16" '0.0000174'
$ws.Range("E16").Value = '  -2.78%  '
Set-TextValue "D17" '24.76'
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("D18").Value = '2.508.90'
$ws.Range("E18").Value = '  -0.93%  '
Set-TextValue "D19" '11.34'
$ws.Range("E19").Value = '  -1.23%  '
Set-TextValue "D20" '7.59'
$ws.Range("E20").Value = '  +0.64%  '
Set-TextValue "D21" '348.29'
$ws.Range("E21").Value = '  -1.85%  '
Set-TextValue "D22" '3.92'
$ws.Range("E22").Value = '  -0.44%  '
$ws.Range("E23").Value = '  +0.76%  '
Set-TextValue "D25" '70.32'
$ws.Range("E25").Value = '  +1.58%  '
Set-TextValue "D26" '3.97'
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("E27").Value = '  -2.59%  '
$ws.Range("D28").Value = '2.649.23'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("D30").Value = '0.0₃0890'
$ws.Range("E30").Value = '  -2.19%  '
Set-TextValue "D31" '7.83'
$ws.Range("E31").Value = '  -0.23%  '
Set-TextValue "D32" '461.55'
$ws.Range("E32").Value = '  -3.55%  '
Set-TextValue "D33" '1.23'
$ws.Range("E33").Value = '  -3.96%  '
$ws.Range("E34").Value = '  -2.01%  '
Set-TextValue "D35" '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +1.70%  '
Set-TextValue "D37" '158.13'
$ws.Range("E37").Value = '  +0.48%  '
Set-TextValue "D38" '18.98'
$ws.Range("E38").Value = '  +0.71%  '
Set-TextValue "D39" '18.51'
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("E40").Value = '  -0.06%  '
Set-TextValue "D41" '4.73'
$ws.Range("E41").Value = '  +0.57%  '
Set-TextValue "D42" '0.318'
$ws.Range("E42").Value = '  -0.40%  '
Set-TextValue "D43" '1.60'
$ws.Range("E43").Value = '  -2.69%  '
Set-TextValue "D44" '38.05'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("E45").Value = '  -13.66%  '
Set-TextValue "D46" '2.25'
$ws.Range("E46").Value = '  -5.85%  '
Set-TextValue "D47" '141.63'
$ws.Range("E47").Value = '  -0.01%  '
Set-TextValue "D48" '0.525'
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("E49").Value = '  -1.85%  '
Set-TextValue "D50" '0.0729'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("E51").Value = '  -3.20%  '
